$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source row (22, "Bibliografia:" / long bibliography text) is removed
# entirely - delete it first (from the bottom) so row numbering above is
# untouched while we rewrite rows 13-21 in place.
$ws.Rows.Item(22).Delete()

# Row 13: "Programa resumido:" / "Semestral"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# Row 14: "Short syllabus:" (B/C empty)
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = ""
$ws.Rows.Item(14).RowHeight = 60

# Row 15: "Programa:" / "01/01/2018"
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"
$ws.Rows.Item(15).RowHeight = 120

# Row 16: "Syllabus:" (B/C empty)
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = ""
$ws.Rows.Item(16).RowHeight = 120

# Row 17: "Avaliação:" (B/C empty, default row height)
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = ""
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(17).Insert()
$ws.Range("A17").Value = "Avaliação:"

# Row 18: "Método:" / "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C18").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Rows.Item(18).RowHeight = 60

# Row 19: "Critério:" / "Atividades docentes: ..."
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Atividades docentes: Mentoria, palestras e seminários.Atividades discentes: Elaboração de projeto utilizando laboratórios e instalações da USP."
$ws.Range("C19").Value = "Atividades docentes: Mentoria, palestras e seminários.Atividades discentes: Elaboração de projeto utilizando laboratórios e instalações da USP."
$ws.Rows.Item(19).RowHeight = 60

# Row 20: "Norma de recuperação:" / "Avaliação pela equipe de mentores, ..."
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Avaliação pela equipe de mentores, considerando critérios, tais como: qualidade técnica da proposta, grau de inovação, viabilidade técnica, dentre outros.Nota de projeto maior ou igual a 5,0 (cinco)."
$ws.Range("C20").Value = "Avaliação pela equipe de mentores, considerando critérios, tais como: qualidade técnica da proposta, grau de inovação, viabilidade técnica, dentre outros.Nota de projeto maior ou igual a 5,0 (cinco)."
$ws.Rows.Item(20).RowHeight = 60

# Row 21: "Bibliografia:" / "Não há recuperação."
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Não há recuperação."
$ws.Range("C21").Value = "Não há recuperação."
$ws.Rows.Item(21).RowHeight = 120
